$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 216-217, shifting existing rows 216.. down to 218..
$ws.Range("A216:A217").EntireRow.Insert()

# Populate new row 216
$ws.Range("A216").Value = 3
$ws.Range("B216").Value = "Femacal de La Calera"
$ws.Range("C216").Value = "Coquimbo"
$ws.Range("D216").Value = 44889
$ws.Range("E216").Value = 5
$ws.Range("F216").Value = "Fruta"
$ws.Range("G216").Value = 100101
$ws.Range("H216").Value = "Berries"
$ws.Range("I216").Value = 100101001
$ws.Range("J216").Value = "Arándano (blue)"
$ws.Range("K216").Value = "Sin especificar"
$ws.Range("L216").Value = "Primera"
$ws.Range("M216").Value = 110
$ws.Range("N216").Value = 5800
$ws.Range("O216").Value = 6000
$ws.Range("P216").Value = 5918
$ws.Range("Q216").Value = "$/bandeja 2 kilos"
$ws.Range("R216").Value = "Provincia de Linares"
$ws.Range("S216").Value = 2959
$ws.Range("T216").Value = 2

# Populate new row 217
$ws.Range("A217").Value = 3
$ws.Range("B217").Value = "Femacal de La Calera"
$ws.Range("C217").Value = "Coquimbo"
$ws.Range("D217").Value = 44889
$ws.Range("E217").Value = 5
$ws.Range("F217").Value = "Fruta"
$ws.Range("G217").Value = 100101
$ws.Range("H217").Value = "Berries"
$ws.Range("I217").Value = 100101001
$ws.Range("J217").Value = "Arándano (blue)"
$ws.Range("K217").Value = "Sin especificar"
$ws.Range("L217").Value = "Segunda"
$ws.Range("M217").Value = 60
$ws.Range("N217").Value = 4000
$ws.Range("O217").Value = 4000
$ws.Range("P217").Value = 4000
$ws.Range("Q217").Value = "$/bandeja 2 kilos"
$ws.Range("R217").Value = "Provincia de Linares"
$ws.Range("S217").Value = 2000
$ws.Range("T217").Value = 2
